$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Time Table")

$ws.Range("A4").Value = "CS291[GY, AH]  /  M201(T)[ABj]"
$ws.Range("B4").Value = "CS291[GY, AH]  /  M201(T)[ABj]"
$ws.Range("C4").Value = "CS291[GY, AH]  /  M201(T)[ABj]"
$ws.Range("D4").Value = "CS291[GY, AH]  /  M201(T)[ABj]"
$ws.Range("E4").Value = "CS201[GY]  /  "
$ws.Range("F4").Value = "M201[ABj]  /  "
$ws.Range("G4").Value = "CH201[SC, PD]  /  "
$ws.Range("A6").Value = "CS291[GY, AH]  /  M201(T)[ABj]"
$ws.Range("B6").Value = "CS291[GY, AH]  /  M201(T)[ABj]"
$ws.Range("C6").Value = "CS291[GY, AH]  /  M201(T)[ABj]"
$ws.Range("D6").Value = "CS291[GY, AH]  /  M201(T)[ABj]"
$ws.Range("E6").Value = "M201[ABj]  /  "
$ws.Range("F6").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("G6").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("A8").Value = "ME291[TR, TKG]  /  "
$ws.Range("B8").Value = "ME291[TR, TKG]  /  "
$ws.Range("C8").Value = "ME291[TR, TKG]  /  "
$ws.Range("D8").Value = "ME291[TR, TKG]  /  "
$ws.Range("E8").Value = "ME291[TR, TKG]  /  "
$ws.Range("F8").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("G8").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("A10").Value = "Free Period!"
$ws.Range("B10").Value = "HU201[SDa]  /  "
$ws.Range("C10").Value = "Free Period!"
$ws.Range("D10").Value = "CH201[SC, PD]  /  "
$ws.Range("F10").Value = "CS201[GY]  /  "
$ws.Range("A12").Value = "M201[ABj]  /  "
$ws.Range("B12").Value = "HU201[SDa]  /  "
$ws.Range("C12").Value = "Free Period!"
$ws.Range("D12").Value = "CH201[SC, PD]  /  "
$ws.Range("E12").Value = "Free Period!"
$ws.Range("F12").Value = "CS201[GY]  /  "
$ws.Range("G12").Value = "Free Period!"
$ws.Range("A16").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("B16").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("C16").Value = "Free Period!"
$ws.Range("E16").Value = "HU201[SDa]  /  "
$ws.Range("F16").Value = "CH201[PD]  /  "
$ws.Range("G16").Value = "Free Period!"
$ws.Range("A18").Value = "ME291[TR, TKG]  /  "
$ws.Range("B18").Value = "ME291[TR, TKG]  /  "
$ws.Range("C18").Value = "ME291[TR, TKG]  /  "
$ws.Range("D18").Value = "ME291[TR, TKG]  /  "
$ws.Range("E18").Value = "ME291[TR, TKG]  /  "
$ws.Range("F18").Value = "M201[ABj]  /  "
$ws.Range("G18").Value = "CH201[PD]  /  "
$ws.Range("A20").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("B20").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("C20").Value = "IT201[AKS]  /  "
$ws.Range("D20").Value = "Free Period!"
$ws.Range("E20").Value = "M201[ABj]  /  "
$ws.Range("F20").Value = "CH201[PD]  /  "
$ws.Range("A22").Value = "IT201[AKS]  /  "
$ws.Range("B22").Value = "M201[ABj]  /  "
$ws.Range("C22").Value = "HU201[SDa]  /  "
$ws.Range("D22").Value = "IT291[AKS, ARC]  /  M201(T)[SCh]"
$ws.Range("E22").Value = "IT291[AKS, ARC]  /  M201(T)[SCh]"
$ws.Range("F22").Value = "IT291[AKS, ARC]  /  M201(T)[SCh]"
$ws.Range("G22").Value = "IT291[AKS, ARC]  /  M201(T)[SCh]"
$ws.Range("A24").Value = "IT201[AKS]  /  "
$ws.Range("B24").Value = "Free Period!"
$ws.Range("C24").Value = "Free Period!"
$ws.Range("D24").Value = "IT291[AKS, ARC]  /  M201(T)[SCh]"
$ws.Range("E24").Value = "IT291[AKS, ARC]  /  M201(T)[SCh]"
$ws.Range("F24").Value = "IT291[AKS, ARC]  /  M201(T)[SCh]"
$ws.Range("G24").Value = "IT291[AKS, ARC]  /  M201(T)[SCh]"
$ws.Range("A28").Value = "ECE291[SDe, SSK]  /  M201(T)[SCh]"
$ws.Range("B28").Value = "ECE291[SDe, SSK]  /  M201(T)[SCh]"
$ws.Range("C28").Value = "ECE291[SDe, SSK]  /  M201(T)[SCh]"
$ws.Range("D28").Value = "ECE291[SDe, SSK]  /  M201(T)[SCh]"
$ws.Range("E28").Value = "ECE201[SMa]  /  "
$ws.Range("F28").Value = "M201[SCh]  /  "
$ws.Range("G28").Value = "PH201[SoM]  /  "
$ws.Range("A30").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("B30").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("C30").Value = "PH201(T)[AT, SoM]  /  "
$ws.Range("D30").Value = "HU201[SDa]  /  "
$ws.Range("E30").Value = "Free Period!"
$ws.Range("G30").Value = "PH201[SoM]  /  "
$ws.Range("A32").Value = "Free Period!"
$ws.Range("B32").Value = "PH201[SoM]  /  "
$ws.Range("C32").Value = "HU201[SDa]  /  "
$ws.Range("D32").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("E32").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("F32").Value = "M201[SCh]  /  "
$ws.Range("G32").Value = "ECE201[SMa]  /  "
$ws.Range("A34").Value = "ECE291[SDe, SSK]  /  M201(T)[SCh]"
$ws.Range("B34").Value = "ECE291[SDe, SSK]  /  M201(T)[SCh]"
$ws.Range("C34").Value = "ECE291[SDe, SSK]  /  M201(T)[SCh]"
$ws.Range("D34").Value = "ECE291[SDe, SSK]  /  M201(T)[SCh]"
$ws.Range("E34").Value = "Free Period!"
$ws.Range("F34").Value = "Free Period!"
$ws.Range("G34").Value = "Free Period!"
$ws.Range("A36").Value = "M201[SCh]  /  "
$ws.Range("B36").Value = "ECE201[SMa]  /  "
$ws.Range("C36").Value = "ME291[BDC]  /  "
$ws.Range("D36").Value = "ME291[BDC]  /  "
$ws.Range("E36").Value = "ME291[BDC]  /  "
$ws.Range("F36").Value = "ME291[BDC]  /  "
$ws.Range("G36").Value = "ME291[BDC]  /  "
$ws.Range("A40").Value = "Free Period!"
$ws.Range("B40").Value = "Free Period!"
$ws.Range("C40").Value = "M201[SCh]  /  "
$ws.Range("D40").Value = "HU201[SDa]  /  "
$ws.Range("E40").Value = "PH201[AT]  /  "
$ws.Range("F40").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("G40").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("A42").Value = "M201(T)[SCh]  /  "
$ws.Range("B42").Value = "EE201[SL]  /  "
$ws.Range("C42").Value = "M201[SCh]  /  "
$ws.Range("D42").Value = "EE291[AKS, SL]  /  PH201(T)[AT, SoM]"
$ws.Range("E42").Value = "EE291[AKS, SL]  /  PH201(T)[AT, SoM]"
$ws.Range("F42").Value = "EE291[AKS, SL]  /  PH201(T)[AT, SoM]"
$ws.Range("G42").Value = "EE291[AKS, SL]  /  PH201(T)[AT, SoM]"
$ws.Range("B44").Value = "EE201[SL]  /  "
$ws.Range("C44").Value = "PH201[AT]  /  "
$ws.Range("D44").Value = "EE291[AKS, SL]  /  PH201(T)[AT, SoM]"
$ws.Range("E44").Value = "EE291[AKS, SL]  /  PH201(T)[AT, SoM]"
$ws.Range("F44").Value = "EE291[AKS, SL]  /  PH201(T)[AT, SoM]"
$ws.Range("G44").Value = "EE291[AKS, SL]  /  PH201(T)[AT, SoM]"
$ws.Range("A46").Value = "ME291[BDC]  /  "
$ws.Range("B46").Value = "ME291[BDC]  /  "
$ws.Range("C46").Value = "ME291[BDC]  /  "
$ws.Range("D46").Value = "ME291[BDC]  /  "
$ws.Range("E46").Value = "ME291[BDC]  /  "
$ws.Range("F46").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("G46").Value = "HU291[SDa]  /  MOOCS[]"
$ws.Range("A48").Value = "Free Period!"
$ws.Range("B48").Value = "Free Period!"
$ws.Range("C48").Value = "PH201[AT]  /  "
$ws.Range("D48").Value = "EE201[SL]  /  "
$ws.Range("E48").Value = "Free Period!"
$ws.Range("F48").Value = "M201[SCh]  /  "
$ws.Range("G48").Value = "HU201[SDa]  /  "
